$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(5, 8).Value = 150
$ws.Cells.Item(5, 9).Value = 150
$ws.Cells.Item(5, 11).Value = 150
$ws.Cells.Item(5, 13).Value = -35
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 6533
$ws.Cells.Item(116, 9).Value = 8666
$ws.Cells.Item(116, 11).Value = 8666
$ws.Cells.Item(116, 13).Value = -5224
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 333.5
$ws.Cells.Item(5, 9).Value = 129.5
$ws.Cells.Item(5, 10).Value = 537.5
$ws.Cells.Item(5, 11).Value = 129.5
$ws.Cells.Item(5, 12).Value = 537.5
$ws.Cells.Item(5, 13).Value = -17.5
$ws.Cells.Item(5, 14).Value = -761.5
$ws.Cells.Item(31, 8).Value = 8748.200000000001
$ws.Cells.Item(31, 9).Value = 8748.200000000001
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 8748.200000000001
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -8454.200000000001
$ws.Cells.Item(31, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2212.375
$ws.Cells.Item(122, 9).Value = 1928.6154
$ws.Cells.Item(122, 10).Value = 3442
$ws.Cells.Item(122, 11).Value = 5785.8462
$ws.Cells.Item(122, 12).Value = 10326
$ws.Cells.Item(122, 13).Value = -3335.8462
$ws.Cells.Item(122, 14).Value = -15226
$ws.Cells.Item(124, 8).Value = 46809.668
$ws.Cells.Item(124, 10).Value = 46809.668
$ws.Cells.Item(124, 12).Value = 46809.668
$ws.Cells.Item(124, 14).Value = -56629.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 333.5
$ws.Cells.Item(4, 9).Value = 129.5
$ws.Cells.Item(4, 10).Value = 537.5
$ws.Cells.Item(4, 11).Value = 129.5
$ws.Cells.Item(4, 12).Value = 537.5
$ws.Cells.Item(4, 13).Value = -14.5
$ws.Cells.Item(4, 14).Value = -767.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 14302.75
$ws.Cells.Item(19, 9).Value = 68.333336
$ws.Cells.Item(19, 11).Value = 68.333336
$ws.Cells.Item(19, 13).Value = 101.666664
$ws.Cells.Item(24, 8).Value = 14302.75
$ws.Cells.Item(24, 9).Value = 68.333336
$ws.Cells.Item(24, 11).Value = 68.333336
$ws.Cells.Item(24, 13).Value = 101.666664
$ws.Cells.Item(53, 8).Value = 52092
$ws.Cells.Item(53, 10).Value = 52092
$ws.Cells.Item(53, 12).Value = 52092
$ws.Cells.Item(53, 14).Value = -53306
$ws.Cells.Item(59, 8).Value = 25470.625
$ws.Cells.Item(59, 10).Value = 42691.25
$ws.Cells.Item(59, 12).Value = 42691.25
$ws.Cells.Item(59, 14).Value = -44981.25
$ws.Cells.Item(64, 8).Value = 47750
$ws.Cells.Item(64, 10).Value = 47750
$ws.Cells.Item(64, 12).Value = 47750
$ws.Cells.Item(64, 14).Value = -48246
$ws.Cells.Item(67, 8).Value = 47750
$ws.Cells.Item(67, 10).Value = 47750
$ws.Cells.Item(67, 12).Value = 47750
$ws.Cells.Item(67, 14).Value = -49466
$ws.Cells.Item(99, 8).Value = 2634.4546
$ws.Cells.Item(99, 9).Value = 2122.25
$ws.Cells.Item(99, 11).Value = 2122.25
$ws.Cells.Item(99, 13).Value = -624.25
$ws.Cells.Item(126, 8).Value = 2634.4546
$ws.Cells.Item(126, 9).Value = 2122.25
$ws.Cells.Item(126, 11).Value = 6366.75
$ws.Cells.Item(126, 13).Value = -3896.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 133551.8
$ws.Cells.Item(4, 10).Value = 282.7
$ws.Cells.Item(4, 12).Value = 848.0999999999999
$ws.Cells.Item(4, 14).Value = -1072.1
$ws.Cells.Item(11, 8).Value = 50119.35
$ws.Cells.Item(11, 9).Value = 55671.5
$ws.Cells.Item(11, 10).Value = 150
$ws.Cells.Item(11, 11).Value = 167014.5
$ws.Cells.Item(11, 12).Value = 450
$ws.Cells.Item(11, 13).Value = -166874.5
$ws.Cells.Item(11, 14).Value = -730
$ws.Cells.Item(63, 8).Value = 1304.6666
$ws.Cells.Item(63, 9).Value = 82
$ws.Cells.Item(63, 10).Value = 3750
$ws.Cells.Item(63, 11).Value = 246
$ws.Cells.Item(63, 12).Value = 11250
$ws.Cells.Item(63, 13).Value = 503
$ws.Cells.Item(63, 14).Value = -12748
$ws.Cells.Item(66, 8).Value = 1304.6666
$ws.Cells.Item(66, 9).Value = 82
$ws.Cells.Item(66, 10).Value = 3750
$ws.Cells.Item(66, 11).Value = 738
$ws.Cells.Item(66, 12).Value = 33750
$ws.Cells.Item(66, 13).Value = 3006
$ws.Cells.Item(66, 14).Value = -41238
$ws.Cells.Item(99, 8).Value = 50000
$ws.Cells.Item(99, 9).Value = 50000
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 150000
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -147754
$ws.Cells.Item(99, 14).ClearContents()
$ws.Cells.Item(108, 8).Value = 837
$ws.Cells.Item(108, 9).Value = 398.33334
$ws.Cells.Item(108, 11).Value = 1195.00002
$ws.Cells.Item(108, 13).Value = 1684.99998
$ws.Cells.Item(126, 8).Value = 2000
$ws.Cells.Item(126, 9).Value = 2000
$ws.Cells.Item(126, 11).Value = 6000
$ws.Cells.Item(126, 13).Value = -1060
$ws.Cells.Item(132, 8).Value = 632.8889
$ws.Cells.Item(132, 9).Value = 399.42856
$ws.Cells.Item(132, 10).Value = 1450
$ws.Cells.Item(132, 11).Value = 3594.85704
$ws.Cells.Item(132, 12).Value = 13050
$ws.Cells.Item(132, 13).Value = -1064.85704
$ws.Cells.Item(132, 14).Value = -18110
$ws.Cells.Item(134, 8).Value = 3961
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 3349.125
$ws.Cells.Item(139, 9).Value = 1940.5
$ws.Cells.Item(139, 10).Value = 7575
$ws.Cells.Item(139, 11).Value = 5821.5
$ws.Cells.Item(139, 12).Value = 22725
$ws.Cells.Item(139, 13).Value = -681.5
$ws.Cells.Item(139, 14).Value = -33005
$ws.Cells.Item(140, 8).Value = 1654.5
$ws.Cells.Item(140, 9).Value = 1397.1538
$ws.Cells.Item(140, 11).Value = 4191.4614
$ws.Cells.Item(140, 13).Value = 988.5385999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 42510
$ws.Cells.Item(52, 10).Value = 55000
$ws.Cells.Item(52, 12).Value = 55000
$ws.Cells.Item(52, 14).Value = -55518
$ws.Cells.Item(122, 8).Value = 1744.625
$ws.Cells.Item(122, 9).Value = 1149.75
$ws.Cells.Item(122, 10).Value = 2339.5
$ws.Cells.Item(122, 11).Value = 3449.25
$ws.Cells.Item(122, 12).Value = 7018.5
$ws.Cells.Item(122, 13).Value = -999.25
$ws.Cells.Item(122, 14).Value = -11918.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 360.75
$ws.Cells.Item(16, 9).Value = 215.11111
$ws.Cells.Item(16, 11).Value = 215.11111
$ws.Cells.Item(16, 13).Value = -45.11111
$ws.Cells.Item(63, 8).Value = 44444
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 8).Value = 44444
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 5957.909
$ws.Cells.Item(68, 9).Value = 4187.4
$ws.Cells.Item(68, 11).Value = 4187.4
$ws.Cells.Item(68, 13).Value = -3438.4
$ws.Cells.Item(71, 8).Value = 5957.909
$ws.Cells.Item(71, 9).Value = 4187.4
$ws.Cells.Item(71, 11).Value = 20937
$ws.Cells.Item(71, 13).Value = -17193
$ws.Cells.Item(127, 8).Value = 70001
$ws.Cells.Item(127, 10).Value = 70001
$ws.Cells.Item(127, 12).Value = 70001
$ws.Cells.Item(127, 14).Value = -79921
$ws.Cells.Item(134, 8).Value = 87714.5
$ws.Cells.Item(134, 10).Value = 87714.5
$ws.Cells.Item(134, 12).Value = 87714.5
$ws.Cells.Item(134, 14).Value = -97854.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 49999
$ws.Cells.Item(119, 10).Value = 49999
$ws.Cells.Item(119, 12).Value = 49999
$ws.Cells.Item(119, 14).Value = -59675
$ws.Cells.Item(126, 8).Value = 3273.08
$ws.Cells.Item(126, 9).Value = 1389.25
$ws.Cells.Item(126, 11).Value = 4167.75
$ws.Cells.Item(126, 13).Value = -1697.75
